# Refactor "Pays" (country) column: replace full country names with
# their ISO-3 codes throughout the shared-strings table. Since every
# cell in column D points at one of these 6 shared strings, rewriting
# the string text updates every row at once.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$countryMap = @{
    "Cameroun"            = "CMR"
    "Congo"               = "CNG"
    "Gabon"               = "GAB"
    "Guinée Equatoriale"  = "GNQ"
    "RCA"                 = "CAF"
    "Tchad"               = "TCD"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 76 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null -and $countryMap.ContainsKey($val)) {
        $cell.Value = $countryMap[$val]
    }
}

# The two BGFI / CCEIBANK Guinée-Equatoriale rows no longer need their
# old wrapped-text row height now that the sheet content was touched -
# autofit to let Excel recompute the natural row height.
$ws.Rows.Item(54).AutoFit()
$ws.Rows.Item(55).AutoFit()

# Move the cursor down to the bottom of the list (where the edits were
# made) and select D75.
$ws.Range("D75").Select()
